$d = $word.ActiveDocument

function Get-ParagraphEnd($doc, $pos) {
    $paras = $doc.Paragraphs
    $result = -1
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        $s = $p.Range.Start
        $e = $p.Range.End
        if ($pos -ge $s -and $pos -lt $e) {
            $result = $e
        }
    }
    return $result
}

# ---------------------------------------------------------------
# Change 1: "Review the sales data for ..." bullet paragraph.
# Replace from the start of "last 3 month of previous year" through
# the end of the paragraph text (before the pilcrow) with three runs:
#   "last 3 month, so that we get the trending product for "
#   "up"
#   "coming sale."
# ---------------------------------------------------------------
$f1 = $d.Content
$f1.Find.Execute("last 3 month of previous year", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start1 = $f1.Start
$paraEnd1 = Get-ParagraphEnd $d $start1
$target1 = $d.Range($start1, $paraEnd1 - 1)

$rpr = '<w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
    '<w:r>' + $rpr + '<w:t xml:space="preserve">last 3 month, so that we get the trending product for </w:t></w:r>' + `
    '<w:r>' + $rpr + '<w:t>up</w:t></w:r>' + `
    '<w:r>' + $rpr + '<w:t>coming sale.</w:t></w:r>' + `
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target1.InsertXML($xml1)

# ---------------------------------------------------------------
# Change 2: "Finding out most sold category ..." bullet paragraph.
# Replace from the start of "Finding out most sold category" through
# the end of the paragraph text (before the pilcrow) with the new,
# more-granular run/proofErr structure.
# ---------------------------------------------------------------
$f2 = $d.Content
$f2.Find.Execute("Finding out most sold category", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start2 = $f2.Start
$paraEnd2 = Get-ParagraphEnd $d $start2
$target2 = $d.Range($start2, $paraEnd2 - 1)

$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
    '<w:r>' + $rpr + '<w:t xml:space="preserve">Finding out most sold category from previous </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r>' + $rpr + '<w:t>3 month</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r>' + $rpr + '<w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r>' + $rpr + '<w:t xml:space="preserve">sales data </w:t></w:r>' + `
    '<w:r>' + $rpr + '<w:t xml:space="preserve">and include </w:t></w:r>' + `
    '<w:r>' + $rpr + '<w:t>all product</w:t></w:r>' + `
    '<w:r>' + $rpr + '<w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r>' + $rpr + '<w:t>within this</w:t></w:r>' + `
    '<w:r>' + $rpr + '<w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r>' + $rpr + '<w:t>category in</w:t></w:r>' + `
    '<w:r>' + $rpr + '<w:t xml:space="preserve"> Big day sale except </w:t></w:r>' + `
    '<w:r>' + $rpr + '<w:t>backorders.</w:t></w:r>' + `
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target2.InsertXML($xml2)

# ---------------------------------------------------------------
# Change 3: move the _GoBack bookmark from the last (empty) paragraph
# to the second-to-last (empty) paragraph.
# ---------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count
$secondLastPara = $paras.Item($count - 1)

foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

$d.Bookmarks.Add("_GoBack", $secondLastPara.Range)
